{"js": "// Office.js (Word JavaScript API) script implementing the commit's text\n// edits:\n//  1. \"...and provide processed data to different environments.\" -- no text\n//     change, but the \"_GoBack\" bookmark that used to sit in front of the\n//     \"Work History\" heading is now placed right before \"processed data to \"\n//     (splitting that run into \"and provide \" + \"processed data to \").\n//  2. The \"In-depth knowledge of Hadoop...\" bullet gets reworded to mention\n//     Spark architectures and adds Executor/Driver/Cluster Manager.\n//  3. The \"Extensive knowledge in reviewing Hadoop log files...\" bullet gets\n//     reworded (\"code modules / applications\") and its tooling list is\n//     rephrased.\n\nconst body = context.document.body;\n\n// --- 1. Relocate the \"_GoBack\" bookmark ------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nlet results = body.search(\"processed data to\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst goBackStart = results.items[0].getRange(\"Start\");\ngoBackStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Reword the \"In-depth knowledge...\" bullet --------------------------\nresults = body.search(\n  \"In-depth knowledge of Hadoop architecture and various components such as HDFS, Resource Manager, Application Master, Node Manager, Name Node, Data Node and Map-Reduce concepts\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\n  \"In-depth knowledge of Hadoop, Spark architectures and various components such as HDFS, Resource Manager, Application Master, Node Manager, Name Node, Data Node and Map-Reduce concepts, Executor, Driver, Cluster Manager\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 3. Reword the \"Extensive knowledge...\" bullet -------------------------\nresults = body.search(\"developing code using hive\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\n  \"developing code modules / applications using hive\",\n  \"Replace\"\n);\nawait context.sync();\n\nresults = body.search(\n  \", Pig Latin, Impala and Hadoop Map-Reduce applications, Spark applications in Scala\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\n  \", Pig Latin, Impala, Hadoop, Map-Reduce and Spark\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop script implementing the commit's text edits:\n#  1. \"...and provide processed data to different environments.\" -- no text\n#     change, but the \"_GoBack\" bookmark that used to sit in front of the\n#     \"Work History\" heading is now placed right before \"processed data to \"\n#     (splitting that run into \"and provide \" + \"processed data to \").\n#  2. The \"In-depth knowledge of Hadoop...\" bullet gets reworded to mention\n#     Spark architectures and adds Executor/Driver/Cluster Manager.\n#  3. The \"Extensive knowledge in reviewing Hadoop log files...\" bullet gets\n#     reworded (\"code modules / applications\") and its tooling list is\n#     rephrased.\n\n$d = $word.ActiveDocument\n\n# --- 1. Relocate the \"_GoBack\" bookmark -----------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"processed data to\") | Out-Null\n$rng.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n\n# --- 2. Reword the \"In-depth knowledge...\" bullet -------------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"In-depth knowledge of Hadoop architecture and various components such as HDFS, Resource Manager, Application Master, Node Manager, Name Node, Data Node and Map-Reduce concepts\") | Out-Null\n$rng2.Text = \"In-depth knowledge of Hadoop, Spark architectures and various components such as HDFS, Resource Manager, Application Master, Node Manager, Name Node, Data Node and Map-Reduce concepts, Executor, Driver, Cluster Manager\"\n\n# --- 3. Reword the \"Extensive knowledge...\" bullet -------------------------\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Execute(\"developing code using hive\") | Out-Null\n$rng3.Text = \"developing code modules / applications using hive\"\n\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.Execute(\", Pig Latin, Impala and Hadoop Map-Reduce applications, Spark applications in Scala\") | Out-Null\n$rng4.Text = \", Pig Latin, Impala, Hadoop, Map-Reduce and Spark\"\n"}
